$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.956.96"
$ws.Range("E2").Value = "  -2.20%  "

$ws.Range("D3").Value = "2.255.58"

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "298.34"
$ws.Range("E5").Value = "  -2.70%  "

$ws.Range("D6").Value = "93.37"
$ws.Range("E6").Value = "  -6.86%  "

$ws.Range("E7").Value = "  -2.78%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -3.86%  "

$ws.Range("D10").Value = "32.86"
$ws.Range("E10").Value = "  -5.92%  "

$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("D12").Value = "47.74"
$ws.Range("E12").Value = "  -8.42%  "

$ws.Range("D13").Value = "0.112"
$ws.Range("E13").Value = "  +0.27%  "

$ws.Range("D14").Value = "6.66"

$ws.Range("D15").Value = "2.606.51"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("D16").Value = "15.31"
$ws.Range("E16").Value = "  -3.92%  "

$ws.Range("D17").Value = "2.245.19"
$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.770"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.30%  "

$ws.Range("D19").Value = "41.975.84"
$ws.Range("E19").Value = "  -2.03%  "

$ws.Range("D20").Value = "0.0₃0889"
$ws.Range("E20").Value = "  -2.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.54%  "

$ws.Range("D22").Value = "11.31"
$ws.Range("E22").Value = "  -3.63%  "

$ws.Range("D23").Value = "66.42"
$ws.Range("E23").Value = "  -2.22%  "

$ws.Range("D24").Value = "233.54"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").Value = "1.92"
$ws.Range("E25").Value = "  -5.45%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "2.44"
$ws.Range("E27").Value = "  -4.66%  "

$ws.Range("D28").Value = "23.62"
$ws.Range("E28").Value = "  -7.37%  "

$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -6.62%  "

$ws.Range("D30").Value = "167.68"
$ws.Range("E30").Value = "  +4.68%  "

$ws.Range("D31").Value = "33.46"

$ws.Range("D32").Value = "9.01"
$ws.Range("E32").Value = "  -3.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.09%  "

$ws.Range("D34").Value = "4.92"
$ws.Range("E34").Value = "  -4.06%  "

$ws.Range("E35").Value = "  -4.68%  "

$ws.Range("E36").Value = "  -5.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0690"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.19%  "

$ws.Range("D38").Value = "2.77"
$ws.Range("E38").Value = "  -5.92%  "

$ws.Range("D39").Value = "15.87"
$ws.Range("E39").Value = "  -8.71%  "

$ws.Range("D40").Value = "0.0989"
$ws.Range("E40").Value = "  -3.62%  "

$ws.Range("E41").Value = "  -3.60%  "

$ws.Range("E42").Value = "  -8.77%  "

$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  -1.03%  "

$ws.Range("D44").Value = "1.947.32"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("D45").Value = "0.0277"
$ws.Range("E45").Value = "  -2.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.20"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.21%  "

$ws.Range("D47").Value = "9.53"
$ws.Range("E47").Value = "  -7.42%  "

$ws.Range("D48").Value = "2.77"
$ws.Range("E48").Value = "  -4.91%  "

$ws.Range("E49").Value = "  -2.48%  "

$ws.Range("D50").Value = "2.480.75"
$ws.Range("E50").Value = "  -3.05%  "

$ws.Range("D51").Value = "51.78"
$ws.Range("E51").Value = "  -7.35%  "
